$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new column H: header "Tiền phạt" and value 0 for each of the 20 employee rows
$ws.Range("H1").Value = "Tiền phạt"

for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}

# Update the active selection as per the diff (J21)
$ws.Range("J21").Select()
